$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 152, shifting existing rows 152-190 down to 153-191.
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row 152 with the new record.
$ws.Cells.Item(152, 1).Value = 5
$ws.Cells.Item(152, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(152, 3).Value = "Maule"
$ws.Cells.Item(152, 4).Value = 45135
$ws.Cells.Item(152, 5).Value = 7
$ws.Cells.Item(152, 6).Value = 100112001
$ws.Cells.Item(152, 7).Value = "Berenjena"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 200
$ws.Cells.Item(152, 11).Value = 7000
$ws.Cells.Item(152, 12).Value = 7000
$ws.Cells.Item(152, 13).Value = 7000
$ws.Cells.Item(152, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(152, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(152, 16).Value = 140
$ws.Cells.Item(152, 17).Value = 50
$ws.Cells.Item(152, 18).Value = "Hortaliza"
